$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Materialise the previously-blank cells on existing rows (2 & 3) ---
# Copy the format of an already-blank, already-styled cell (C2) onto the
# new blank cells so they pick up the very same style entry instead of
# minting a new (and possibly orphaned) one.
$ws.Range("C2").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("C3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 3: Tuesday (C3) is now also marked available.
$ws.Range("C3").Value = $true

# --- New row 4: a new lecturer, available only on Thursday ---
$ws.Range("A4").Value = 57381920
$ws.Range("E4").Value = $true

# Apply the new wrap-text / vertically-centered style to row 4.
# Set it once on A4 then propagate via a format-only paste so every
# cell in the row shares the very same style entry.
$a4 = $ws.Range("A4")
$a4.WrapText = $true
$a4.VerticalAlignment = -4108

$a4.Copy()
$ws.Range("B4:F4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move the selection to the newly added row, like the source workbook.
$ws.Range("A4").Select()
